$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "71.977.68"
$ws.Range("E2").Value = "  -0.97%  "

$ws.Range("D3").Value = "2.664.56"
$ws.Range("E3").Value = "  +1.41%  "

$ws.Range("E4").Value = "  -0.03%  "

$ws.Range("D5").Value = "'598.61"
$ws.Range("E5").Value = "  -1.00%  "

$ws.Range("D6").Value = "'174.63"
$ws.Range("E6").Value = "  -2.14%  "

$ws.Range("E7").Value = "  -0.02%  "

$ws.Range("D8").Value = "'0.524"
$ws.Range("E8").Value = "  -0.47%  "

$ws.Range("D9").Value = "2.662.64"
$ws.Range("E9").Value = "  +1.41%  "

$ws.Range("E10").Value = "  -1.33%  "

$ws.Range("E11").Value = "  +2.43%  "

$ws.Range("E12").Value = "  +0.83%  "

$ws.Range("E13").Value = "  -1.24%  "

$ws.Range("D14").Value = "3.149.25"
$ws.Range("E14").Value = "  +0.35%  "

$ws.Range("D15").Value = "'0.0000184"
$ws.Range("E15").Value = "  -2.42%  "

$ws.Range("D16").Value = "71.845.15"
$ws.Range("E16").Value = "  -1.10%  "

$ws.Range("D17").Value = "'26.25"

$ws.Range("D18").Value = "2.667.79"
$ws.Range("E18").Value = "  +1.26%  "

$ws.Range("D19").Value = "'12.27"
$ws.Range("E19").Value = "  +6.39%  "

$ws.Range("E20").Value = "  +2.02%  "

$ws.Range("D21").Value = "'371.27"
$ws.Range("E21").Value = "  -3.42%  "

$ws.Range("E22").Value = "  -0.20%  "

$ws.Range("E23").Value = "  +1.42%  "

$ws.Range("D24").Value = "'72.07"
$ws.Range("E24").Value = "  -1.23%  "

$ws.Range("E25").Value = "  -0.05%  "

$ws.Range("E26").Value = "  -1.21%  "

$ws.Range("D27").Value = "'9.74"
$ws.Range("E27").Value = "  -1.13%  "

$ws.Range("D28").Value = "2.801.87"
$ws.Range("E28").Value = "  +1.28%  "

$ws.Range("E29").Value = "  -0.03%  "

$ws.Range("D30").Value = "0.0₃0969"

$ws.Range("E31").Value = "  +0.18%  "

$ws.Range("D32").Value = "'500.57"
$ws.Range("E32").Value = "  -5.70%  "

$ws.Range("E33").Value = "  -1.84%  "

$ws.Range("E34").Value = "  -0.30%  "

$ws.Range("E35").Value = "  -0.09%  "

$ws.Range("D36").Value = "'163.02"
$ws.Range("E36").Value = "  -0.70%  "

$ws.Range("E37").Value = "  +1.02%  "

$ws.Range("E38").Value = "  -0.34%  "

$ws.Range("B39").Value = "Kaspa"
$ws.Range("C39").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D39").Value = "'0.110"
$ws.Range("E39").Value = "  +0.00%  "

$ws.Range("B40").Value = "ImmutableX"
$ws.Range("C40").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D40").Value = "'1.38"
$ws.Range("E40").Value = "  -1.78%  "

$ws.Range("E41").Value = "  -3.18%  "

$ws.Range("E42").Value = "  +0.03%  "

$ws.Range("E43").Value = "  -1.43%  "

$ws.Range("E44").Value = "  -1.81%  "

$ws.Range("D45").Value = "'0.332"
$ws.Range("E45").Value = "  +0.03%  "

$ws.Range("D46").Value = "'156.23"
$ws.Range("E46").Value = "  +3.51%  "

$ws.Range("D47").Value = "'39.48"
$ws.Range("E47").Value = "  -0.05%  "

$ws.Range("E48").Value = "  +3.33%  "

$ws.Range("E49").Value = "  +0.90%  "

$ws.Range("E50").Value = "  +2.16%  "

$ws.Range("E51").Value = "  -1.62%  "
